$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5: update the title text for the flying-emissions regression ---
$ws.Range("E5").Value = "Regression for emissions from flying: OLS (1) and Quantile Regression at 90th (2), 95th (3) and 99th (4) percentile. "

# --- Row 8 (new): total emissions including business travel ---
$ws.Range("E8").Value = "Regression for total emissions, including business travel: OLS (1) and Quantile Regression at 50th (2), 75th (3) and 90th (4) percentile. "
$ws.Range("A8").Value = "total_emissions_withbusiness"
$ws.Range("B8").Value = "Total_emissions"
$ws.Range("C2").Copy($ws.Range("C8"))

# --- Row 9 (new): long-distance emissions including business travel ---
$ws.Range("B9").Value = "emissions_reise"
$ws.Range("A9").Value = "reisen_emissions_withbusiness"
$ws.Range("E9").Value = "Regression for long-distance emissions, including business travels: OLS (1) and Quantile Regression at 50th (2), 75th (3) and 90th (4) percentile. "
$ws.Range("C2").Copy($ws.Range("C9"))

# --- Row 10 (new): flying emissions including business travel ---
$ws.Range("A10").Value = "plane_emissions_withbusiness"
$ws.Range("B10").Value = "emissions_flugzeug"
$ws.Range("E10").Value = "Regression for emissions from flying, including business travels: OLS (1) and Quantile Regression at 90th (2), 95th (3) and 99th (4) percentile.  "
$ws.Range("C2").Copy($ws.Range("C10"))

# --- Column widths (B, C, D got slightly wider) ---
$ws.Columns.Item(2).ColumnWidth = 24.5
$ws.Columns.Item(3).ColumnWidth = 18
$ws.Columns.Item(4).ColumnWidth = 42.5

# --- Selection moves to the new last cell, E10 ---
[void]$ws.Range("E10").Select()
